$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename Sheet3 -> ExamSubjectWizardData and add a new sheet after it
#    named MarksEntry_ExamSubjectWizard.
# ---------------------------------------------------------------------
$wsExam = $wb.Worksheets.Item("Sheet3")
$wsExam.Name = "ExamSubjectWizardData"

$wsMarks = $wb.Worksheets.Add($null, $wsExam)
$wsMarks.Name = "MarksEntry_ExamSubjectWizard"

$wsCategory = $wb.Worksheets.Item("CategorySubjectMappingData")
$wsGroup    = $wb.Worksheets.Item("CategorySubMapping_GridView")

# ---------------------------------------------------------------------
# 2. Populate ExamSubjectWizardData (previously empty Sheet3) with data,
#    mirroring the look of the first two sheets.
# ---------------------------------------------------------------------
$wsExam.Range("A1").Value = "Academic Year"
$wsExam.Range("B1").Value = "Grade"
$wsExam.Range("C1").Value = "Category"
$wsExam.Range("D1").Value = "runMode"

$wsExam.Range("A1:D1").HorizontalAlignment = $wsCategory.Range("A1").HorizontalAlignment
$wsExam.Range("A1:D1").VerticalAlignment = $wsCategory.Range("A1").VerticalAlignment
$wsExam.Range("A1:D1").Font.Bold = $true

$wsExam.Range("A2").Formula = "'2018-2019"
$wsExam.Range("B2").Value = "UKG - EXAM"
$wsExam.Range("C2").Value = "Pre Nursery"
$wsExam.Range("D2").Value = "Y"

$wsExam.Range("A2").HorizontalAlignment = $wsCategory.Range("A2").HorizontalAlignment
$wsExam.Range("A2").VerticalAlignment = $wsCategory.Range("A2").VerticalAlignment

$wsExam.Range("B2:D2").HorizontalAlignment = $wsCategory.Range("C2").HorizontalAlignment
$wsExam.Range("B2:D2").VerticalAlignment = $wsCategory.Range("C2").VerticalAlignment

$wsExam.Columns.Item(1).ColumnWidth = 14
$wsExam.Columns.Item(2).ColumnWidth = 11.42578125
$wsExam.Columns.Item(3).ColumnWidth = 11.42578125

$wsExam.PageSetup.Orientation = 1

$exWin = $wsExam.Windows.Item(1)
$wsExam.Range("E14").Select()

# ---------------------------------------------------------------------
# 3. Populate MarksEntry_ExamSubjectWizard (the brand new sheet).
# ---------------------------------------------------------------------
$wsMarks.Range("A1").Value = "Max Marks"
$wsMarks.Range("B1").Value = "Min Marks"
$wsMarks.Range("C1").Value = "Max Entry Marks"
$wsMarks.Range("D1").Value = "runMode"

$wsMarks.Range("A1:D1").HorizontalAlignment = $wsCategory.Range("A2").HorizontalAlignment
$wsMarks.Range("A1:D1").VerticalAlignment = $wsCategory.Range("A2").VerticalAlignment

$wsMarks.Range("A2").Formula = "'50"
$wsMarks.Range("B2").Formula = "'15"
$wsMarks.Range("C2").Formula = "'50"
$wsMarks.Range("D2").Value = "Y"

$wsMarks.Range("A2:C2").HorizontalAlignment = $wsCategory.Range("A2").HorizontalAlignment
$wsMarks.Range("A2:C2").VerticalAlignment = $wsCategory.Range("A2").VerticalAlignment

$wsMarks.Range("D2").HorizontalAlignment = $wsCategory.Range("C2").HorizontalAlignment
$wsMarks.Range("D2").VerticalAlignment = $wsCategory.Range("C2").VerticalAlignment

$wsMarks.Columns.Item(1).ColumnWidth = 10.42578125
$wsMarks.Columns.Item(2).ColumnWidth = 10.140625
$wsMarks.Columns.Item(3).ColumnWidth = 15.5703125

$wsMarks.Range("A1:D2").Select()

# ---------------------------------------------------------------------
# 4. Update selections on the first two (already existing) sheets.
# ---------------------------------------------------------------------
$wsGroup.Range("D11:D12").Select()
$wsGroup.Application.ActiveWindow.RangeSelection.Item(1).Activate() | Out-Null

# ---------------------------------------------------------------------
# 5. Activate the new MarksEntry_ExamSubjectWizard sheet and set the
#    workbook's first visible sheet to the 2nd tab, matching the diff.
# ---------------------------------------------------------------------
$wsMarks.Activate()
$excel.ActiveWindow.DisplayedSheets = 1
